$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1324.25
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1265.6666
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1265.6666
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -1615.6666
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H137").Value = 4977.077
$ws.Range("I137").Value = 6044.6665
$ws.Range("J137").Value = 2575
$ws.Range("K137").Value = 18133.9995
$ws.Range("L137").Value = 7725
$ws.Range("M137").Value = -15583.9995
$ws.Range("N137").Value = -12825
$ws.Range("H138").Value = 190062.14
$ws.Range("I138").Value = 11149
$ws.Range("J138").Value = 199821.03
$ws.Range("K138").Value = 33447
$ws.Range("L138").Value = 599463.09
$ws.Range("M138").Value = -28307
$ws.Range("N138").Value = -609743.09

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 414964.84
$ws.Range("I32").Value = 475490.9
$ws.Range("J32").Value = 15492.8
$ws.Range("K32").Value = 475490.9
$ws.Range("L32").Value = 15492.8
$ws.Range("M32").Value = -475203.9
$ws.Range("N32").Value = -16066.8
$ws.Range("H110").Value = 1786.3572
$ws.Range("I110").Value = 1684.0834
$ws.Range("K110").Value = 1684.0834
$ws.Range("M110").Value = 360.9166
$ws.Range("H122").Value = 1575.5807
$ws.Range("I122").Value = 1387.5416
$ws.Range("J122").Value = 2220.2856
$ws.Range("K122").Value = 4162.6248
$ws.Range("L122").Value = 6660.8568
$ws.Range("M122").Value = -1712.6248
$ws.Range("N122").Value = -11560.8568
$ws.Range("H132").Value = 3738
$ws.Range("I132").Value = 2887.5356
$ws.Range("J132").Value = 5569.769
$ws.Range("K132").Value = 8662.606800000001
$ws.Range("L132").Value = 16709.307
$ws.Range("M132").Value = -6132.606800000001
$ws.Range("N132").Value = -21769.307

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5073.5
$ws.Range("I31").Value = 1132.6364
$ws.Range("J31").Value = 9890.111000000001
$ws.Range("K31").Value = 1132.6364
$ws.Range("L31").Value = 9890.111000000001
$ws.Range("M31").Value = -837.6364000000001
$ws.Range("N31").Value = -10480.111
$ws.Range("H34").Value = 5073.5
$ws.Range("I34").Value = 1132.6364
$ws.Range("J34").Value = 9890.111000000001
$ws.Range("K34").Value = 1132.6364
$ws.Range("L34").Value = 9890.111000000001
$ws.Range("M34").Value = -930.6364000000001
$ws.Range("N34").Value = -10294.111
$ws.Range("H122").Value = 1596.9459
$ws.Range("I122").Value = 1168.5
$ws.Range("J122").Value = 1857.7391
$ws.Range("K122").Value = 3505.5
$ws.Range("L122").Value = 5573.2173
$ws.Range("M122").Value = -1055.5
$ws.Range("N122").Value = -10473.2173

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 125000060
$ws.Range("I14").Value = 125000060
$ws.Range("K14").Value = 375000180
$ws.Range("M14").Value = -375000007
$ws.Range("H20").Value = 1458.12
$ws.Range("J20").Value = 1458.12
$ws.Range("L20").Value = 4374.36
$ws.Range("N20").Value = -4828.36
$ws.Range("H49").Value = 7497.5
$ws.Range("J49").Value = 7497.5
$ws.Range("L49").Value = 22492.5
$ws.Range("N49").Value = -22804.5
$ws.Range("H75").Value = 6116
$ws.Range("I75").Value = 3100
$ws.Range("J75").Value = 6870
$ws.Range("K75").Value = 9300
$ws.Range("L75").Value = 20610
$ws.Range("M75").Value = -8302
$ws.Range("N75").Value = -22606
$ws.Range("H78").Value = 6116
$ws.Range("I78").Value = 3100
$ws.Range("J78").Value = 6870
$ws.Range("K78").Value = 27900
$ws.Range("L78").Value = 61830
$ws.Range("M78").Value = -22908
$ws.Range("N78").Value = -71814
$ws.Range("H86").Value = 1376
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1376
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4128
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6500
$ws.Range("H89").Value = 1376
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1376
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 12384
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -24240
$ws.Range("H131").Value = 952.907
$ws.Range("J131").Value = 1019.10254
$ws.Range("L131").Value = 3057.30762
$ws.Range("N131").Value = -13137.30762

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1509.0952
$ws.Range("I97").Value = 1399.4286
$ws.Range("J97").Value = 1728.4286
$ws.Range("K97").Value = 1399.4286
$ws.Range("L97").Value = 1728.4286
$ws.Range("M97").Value = -903.4286
$ws.Range("N97").Value = -2720.4286
$ws.Range("H122").Value = 4348.3237
$ws.Range("I122").Value = 2978.4375
$ws.Range("J122").Value = 5566
$ws.Range("K122").Value = 8935.3125
$ws.Range("L122").Value = 16698
$ws.Range("M122").Value = -6485.3125
$ws.Range("N122").Value = -21598
$ws.Range("H132").Value = 2269.5
$ws.Range("I132").Value = 1773.9445
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 5321.833500000001
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -2791.833500000001
$ws.Range("N132").Value = -18558.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 90913150
$ws.Range("I7").Value = 100003770
$ws.Range("K7").Value = 100003770
$ws.Range("M7").Value = -100003658
$ws.Range("H16").Value = 2275
$ws.Range("I16").Value = 2136.2307
$ws.Range("J16").Value = 2500.5
$ws.Range("K16").Value = 2136.2307
$ws.Range("L16").Value = 2500.5
$ws.Range("M16").Value = -1966.2307
$ws.Range("N16").Value = -2840.5
$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -85990
$ws.Range("H126").Value = 90913150
$ws.Range("I126").Value = 100003770
$ws.Range("K126").Value = 300011310
$ws.Range("M126").Value = -300008840
$ws.Range("H132").Value = 3976.9744
$ws.Range("I132").Value = 3714.6191
$ws.Range("J132").Value = 4283.0557
$ws.Range("K132").Value = 11143.8573
$ws.Range("L132").Value = 12849.1671
$ws.Range("M132").Value = -8613.8573
$ws.Range("N132").Value = -17909.1671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 25429.143
$ws.Range("I14").Value = 25999.666
$ws.Range("J14").Value = 25001.25
$ws.Range("K14").Value = 25999.666
$ws.Range("L14").Value = 25001.25
$ws.Range("M14").Value = -25831.666
$ws.Range("N14").Value = -25337.25
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H112").Value = 35000
$ws.Range("J112").Value = 35000
$ws.Range("L112").Value = 35000
$ws.Range("N112").Value = -37954
$ws.Range("H122").Value = 2307.9
$ws.Range("I122").Value = 2167.5386
$ws.Range("J122").Value = 2568.5715
$ws.Range("K122").Value = 6502.6158
$ws.Range("L122").Value = 7705.7145
$ws.Range("M122").Value = -4052.6158
$ws.Range("N122").Value = -12605.7145
$ws.Range("H126").Value = 1346.3572
$ws.Range("I126").Value = 1276.8
$ws.Range("J126").Value = 1385
$ws.Range("K126").Value = 3830.4
$ws.Range("L126").Value = 4155
$ws.Range("M126").Value = -1360.4
$ws.Range("N126").Value = -9095
